$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1116.5555
$ws.Range("I12").Value = 1777
$ws.Range("J12").Value = 786.3333
$ws.Range("K12").Value = 1777
$ws.Range("L12").Value = 786.3333
$ws.Range("M12").Value = -1607
$ws.Range("N12").Value = -1126.3333
$ws.Range("H17").Value = 1452192.9
$ws.Range("I17").Value = 244.33333
$ws.Range("K17").Value = 732.99999
$ws.Range("M17").Value = -564.99999
$ws.Range("H113").Value = 41670240
$ws.Range("I113").Value = 71431496
$ws.Range("J113").Value = 4478.7
$ws.Range("K113").Value = 71431496
$ws.Range("L113").Value = 4478.7
$ws.Range("M113").Value = -71428242
$ws.Range("N113").Value = -10986.7
$ws.Range("H116").Value = 2848
$ws.Range("I116").Value = 1243.75
$ws.Range("J116").Value = 3917.5
$ws.Range("K116").Value = 1243.75
$ws.Range("L116").Value = 3917.5
$ws.Range("M116").Value = 2198.25
$ws.Range("N116").Value = -10801.5
$ws.Range("H129").Value = 228076.73
$ws.Range("J129").Value = 250846.9
$ws.Range("L129").Value = 752540.7
$ws.Range("N129").Value = -762540.7
$ws.Range("H135").Value = 22736528
$ws.Range("I135").Value = 1185.7142
$ws.Range("J135").Value = 62523376
$ws.Range("K135").Value = 10671.4278
$ws.Range("L135").Value = 562710384
$ws.Range("M135").Value = -8136.427799999999
$ws.Range("N135").Value = -562715454
$ws.Range("H136").Value = 46996
$ws.Range("J136").Value = 46996
$ws.Range("L136").Value = 46996
$ws.Range("N136").Value = -57196
$ws.Range("H141").Value = 2867.5
$ws.Range("I141").Value = 2141.4285
$ws.Range("K141").Value = 6424.2855
$ws.Range("M141").Value = -1244.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2888.4375
$ws.Range("I45").Value = 3785.8572
$ws.Range("J45").Value = 2190.4443
$ws.Range("K45").Value = 3785.8572
$ws.Range("L45").Value = 2190.4443
$ws.Range("M45").Value = -3408.8572
$ws.Range("N45").Value = -2944.4443
$ws.Range("H74").Value = 55556384
$ws.Range("I74").Value = 111111540
$ws.Range("J74").Value = 1229.8889
$ws.Range("K74").Value = 111111540
$ws.Range("L74").Value = 1229.8889
$ws.Range("M74").Value = -111110666
$ws.Range("N74").Value = -2977.8889
$ws.Range("H77").Value = 55556384
$ws.Range("I77").Value = 111111540
$ws.Range("J77").Value = 1229.8889
$ws.Range("K77").Value = 555557700
$ws.Range("L77").Value = 6149.4445
$ws.Range("M77").Value = -555553332
$ws.Range("N77").Value = -14885.4445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3120.524
$ws.Range("I20").Value = 3516.5625
$ws.Range("J20").Value = 1853.2
$ws.Range("K20").Value = 3516.5625
$ws.Range("L20").Value = 1853.2
$ws.Range("M20").Value = -3269.5625
$ws.Range("N20").Value = -2347.2
$ws.Range("H86").Value = 1457.9744
$ws.Range("I86").Value = 1377.3243
$ws.Range("K86").Value = 1377.3243
$ws.Range("M86").Value = -254.3243
$ws.Range("H89").Value = 1457.9744
$ws.Range("I89").Value = 1377.3243
$ws.Range("K89").Value = 6886.6215
$ws.Range("M89").Value = -1270.6215
$ws.Range("H105").Value = 3907
$ws.Range("I105").Value = 4164.375
$ws.Range("J105").Value = 3612.8572
$ws.Range("K105").Value = 4164.375
$ws.Range("L105").Value = 3612.8572
$ws.Range("M105").Value = -2417.375
$ws.Range("N105").Value = -7106.8572
$ws.Range("H107").Value = 1070.2
$ws.Range("I107").Value = 1049.75
$ws.Range("K107").Value = 1049.75
$ws.Range("M107").Value = 870.25
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9605.289000000001
$ws.Range("I31").Value = 14586.542
$ws.Range("J31").Value = 3912.4285
$ws.Range("K31").Value = 14586.542
$ws.Range("L31").Value = 3912.4285
$ws.Range("M31").Value = -14291.542
$ws.Range("N31").Value = -4502.4285
$ws.Range("H34").Value = 9605.289000000001
$ws.Range("I34").Value = 14586.542
$ws.Range("J34").Value = 3912.4285
$ws.Range("K34").Value = 14586.542
$ws.Range("L34").Value = 3912.4285
$ws.Range("M34").Value = -14384.542
$ws.Range("N34").Value = -4316.4285
$ws.Range("H58").Value = 11108.857
$ws.Range("I58").Value = 916.5
$ws.Range("J58").Value = 39333.848
$ws.Range("K58").Value = 916.5
$ws.Range("L58").Value = 39333.848
$ws.Range("M58").Value = -713.5
$ws.Range("N58").Value = -39739.848
$ws.Range("H62").Value = 52634412
$ws.Range("I62").Value = 71431200
$ws.Range("J62").Value = 3402.4
$ws.Range("K62").Value = 71431200
$ws.Range("L62").Value = 3402.4
$ws.Range("M62").Value = -71430576
$ws.Range("N62").Value = -4650.4
$ws.Range("H65").Value = 52634412
$ws.Range("I65").Value = 71431200
$ws.Range("J65").Value = 3402.4
$ws.Range("K65").Value = 357156000
$ws.Range("L65").Value = 17012
$ws.Range("M65").Value = -357152880
$ws.Range("N65").Value = -23252
$ws.Range("H132").Value = 29155.7
$ws.Range("I132").Value = 41692.848
$ws.Range("J132").Value = 5872.4287
$ws.Range("K132").Value = 125078.544
$ws.Range("L132").Value = 17617.2861
$ws.Range("M132").Value = -122548.544
$ws.Range("N132").Value = -22677.2861
$ws.Range("H134").Value = 1293.0714
$ws.Range("I134").Value = 1002.55554
$ws.Range("J134").Value = 1510.9584
$ws.Range("K134").Value = 3007.66662
$ws.Range("L134").Value = 4532.8752
$ws.Range("M134").Value = -472.66662
$ws.Range("N134").Value = -9602.8752
$ws.Range("H136").Value = 11108.857
$ws.Range("I136").Value = 916.5
$ws.Range("J136").Value = 39333.848
$ws.Range("K136").Value = 2749.5
$ws.Range("L136").Value = 118001.544
$ws.Range("M136").Value = -199.5
$ws.Range("N136").Value = -123101.544

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 766.28
$ws.Range("J131").Value = 770.4020400000001
$ws.Range("L131").Value = 2311.20612
$ws.Range("N131").Value = -12391.20612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 70176480
$ws.Range("I122").Value = 27778818
$ws.Range("J122").Value = 142858190
$ws.Range("K122").Value = 83336454
$ws.Range("L122").Value = 428574570
$ws.Range("M122").Value = -83334004
$ws.Range("N122").Value = -428579470
$ws.Range("H132").Value = 51948.773
$ws.Range("I132").Value = 56766.26
$ws.Range("K132").Value = 170298.78
$ws.Range("M132").Value = -167768.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1561.381
$ws.Range("I22").Value = 1295.3077
$ws.Range("J22").Value = 1993.75
$ws.Range("K22").Value = 1295.3077
$ws.Range("L22").Value = 1993.75
$ws.Range("M22").Value = -1000.3077
$ws.Range("N22").Value = -2583.75
$ws.Range("H27").Value = 1561.381
$ws.Range("I27").Value = 1295.3077
$ws.Range("J27").Value = 1993.75
$ws.Range("K27").Value = 1295.3077
$ws.Range("L27").Value = 1993.75
$ws.Range("M27").Value = -1188.3077
$ws.Range("N27").Value = -2207.75
$ws.Range("H40").Value = 3239.5945
$ws.Range("I40").Value = 1880.7693
$ws.Range("J40").Value = 3975.625
$ws.Range("K40").Value = 1880.7693
$ws.Range("L40").Value = 3975.625
$ws.Range("M40").Value = -1744.7693
$ws.Range("N40").Value = -4247.625
$ws.Range("H68").Value = 3066.1667
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497
$ws.Range("H71").Value = 3066.1667
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483
$ws.Range("H122").Value = 936835.9
$ws.Range("I122").Value = 1636737.4
$ws.Range("J122").Value = 3633.889
$ws.Range("K122").Value = 4910212.199999999
$ws.Range("L122").Value = 10901.667
$ws.Range("M122").Value = -4907762.199999999
$ws.Range("N122").Value = -15801.667
$ws.Range("H136").Value = 1882.3529
$ws.Range("I136").Value = 1272.7273
$ws.Range("K136").Value = 3818.1819
$ws.Range("M136").Value = -1268.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4852.4
$ws.Range("I62").Value = 3518.3333
$ws.Range("K62").Value = 3518.3333
$ws.Range("M62").Value = -2894.3333
$ws.Range("H65").Value = 4852.4
$ws.Range("I65").Value = 3518.3333
$ws.Range("K65").Value = 17591.6665
$ws.Range("M65").Value = -14471.6665
$ws.Range("H126").Value = 969.9091
$ws.Range("I126").Value = 843.2222
$ws.Range("J126").Value = 1057.6154
$ws.Range("K126").Value = 2529.6666
$ws.Range("L126").Value = 3172.8462
$ws.Range("M126").Value = -59.66660000000002
$ws.Range("N126").Value = -8112.8462
$ws.Range("H132").Value = 2255.5
$ws.Range("I132").Value = 1449.875
$ws.Range("K132").Value = 4349.625
$ws.Range("M132").Value = -1819.625
$ws.Range("H136").Value = 34484944
$ws.Range("I136").Value = 62502144
$ws.Range("K136").Value = 187506432
$ws.Range("M136").Value = -187503882
